# Refresh market-price-derived Leve profit columns (H:N) across all job sheets
# with the latest scheduled Universalis pull. Values only - no formulas/styles change.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3148
$ws.Range("I38").Value = 166.4
$ws.Range("J38").Value = 6875
$ws.Range("K38").Value = 499.2
$ws.Range("L38").Value = 20625
$ws.Range("M38").Value = -127.2
$ws.Range("N38").Value = -21369
$ws.Range("H53").Value = 1078.1538
$ws.Range("I53").Value = 1280.1
$ws.Range("J53").Value = 405
$ws.Range("K53").Value = 1280.1
$ws.Range("L53").Value = 405
$ws.Range("M53").Value = -643.0999999999999
$ws.Range("N53").Value = -1679

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 299.375
$ws.Range("I5").Value = 237
$ws.Range("J5").Value = 403.33334
$ws.Range("K5").Value = 237
$ws.Range("L5").Value = 403.33334
$ws.Range("M5").Value = -125
$ws.Range("N5").Value = -627.33334
$ws.Range("H32").Value = 8785.741
$ws.Range("I32").Value = 6222.795
$ws.Range("J32").Value = 26959.363
$ws.Range("K32").Value = 6222.795
$ws.Range("L32").Value = 26959.363
$ws.Range("M32").Value = -5935.795
$ws.Range("N32").Value = -27533.363
$ws.Range("H39").Value = 4016
$ws.Range("I39").Value = 4016
$ws.Range("K39").Value = 4016
$ws.Range("M39").Value = -3496
$ws.Range("H74").Value = 2594.3076
$ws.Range("I74").Value = 1768
$ws.Range("J74").Value = 4453.5
$ws.Range("K74").Value = 1768
$ws.Range("L74").Value = 4453.5
$ws.Range("M74").Value = -894
$ws.Range("N74").Value = -6201.5
$ws.Range("H77").Value = 2594.3076
$ws.Range("I77").Value = 1768
$ws.Range("J77").Value = 4453.5
$ws.Range("K77").Value = 8840
$ws.Range("L77").Value = 22267.5
$ws.Range("M77").Value = -4472
$ws.Range("N77").Value = -31003.5
$ws.Range("H97").Value = 583.625
$ws.Range("I97").Value = 569.25
$ws.Range("J97").Value = 655.5
$ws.Range("K97").Value = 569.25
$ws.Range("L97").Value = 655.5
$ws.Range("M97").Value = -73.25
$ws.Range("N97").Value = -1647.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 299.375
$ws.Range("I4").Value = 237
$ws.Range("J4").Value = 403.33334
$ws.Range("K4").Value = 237
$ws.Range("L4").Value = 403.33334
$ws.Range("M4").Value = -122
$ws.Range("N4").Value = -633.33334
$ws.Range("H99").Value = 45455668
$ws.Range("J99").Value = 1240.6666
$ws.Range("L99").Value = 1240.6666
$ws.Range("N99").Value = -4236.6666
$ws.Range("H134").Value = 15577.286
$ws.Range("I134").Value = 1506.8334
$ws.Range("K134").Value = 4520.5002
$ws.Range("M134").Value = -1985.5002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H35").Value = 1112.5
$ws.Range("I35").Value = 1112.5
$ws.Range("K35").Value = 1112.5
$ws.Range("M35").Value = -818.5
$ws.Range("H99").Value = 2056
$ws.Range("I99").Value = 1860
$ws.Range("K99").Value = 1860
$ws.Range("M99").Value = -362
$ws.Range("H107").Value = 1800
$ws.Range("J107").Value = 2300
$ws.Range("L107").Value = 2300
$ws.Range("N107").Value = -6140
$ws.Range("H126").Value = 2056
$ws.Range("I126").Value = 1860
$ws.Range("K126").Value = 5580
$ws.Range("M126").Value = -3110
$ws.Range("H134").Value = 1509.1034
$ws.Range("I134").Value = 1413.4615
$ws.Range("J134").Value = 2338
$ws.Range("K134").Value = 4240.3845
$ws.Range("L134").Value = 7014
$ws.Range("M134").Value = -1705.3845
$ws.Range("N134").Value = -12084

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 16154067
$ws.Range("I131").Value = 71429070
$ws.Range("J131").Value = 32190.041
$ws.Range("K131").Value = 214287210
$ws.Range("L131").Value = 96570.12300000001
$ws.Range("M131").Value = -214282170
$ws.Range("N131").Value = -106650.123

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 2017
$ws.Range("I36").Value = 2017
$ws.Range("K36").Value = 2017
$ws.Range("M36").Value = -1532
$ws.Range("H80").Value = 4888.375
$ws.Range("I80").Value = 4000
$ws.Range("J80").Value = 5015.2856
$ws.Range("K80").Value = 4000
$ws.Range("L80").Value = 5015.2856
$ws.Range("M80").Value = -3002
$ws.Range("N80").Value = -7011.2856
$ws.Range("H83").Value = 4888.375
$ws.Range("I83").Value = 4000
$ws.Range("J83").Value = 5015.2856
$ws.Range("K83").Value = 20000
$ws.Range("L83").Value = 25076.428
$ws.Range("M83").Value = -15008
$ws.Range("N83").Value = -35060.428

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2986.6667
$ws.Range("I7").Value = 2895
$ws.Range("J7").Value = 3170
$ws.Range("K7").Value = 2895
$ws.Range("L7").Value = 3170
$ws.Range("M7").Value = -2783
$ws.Range("N7").Value = -3394
$ws.Range("H40").Value = 4317.1177
$ws.Range("I40").Value = 2215.0833
$ws.Range("J40").Value = 9362
$ws.Range("K40").Value = 2215.0833
$ws.Range("L40").Value = 9362
$ws.Range("M40").Value = -2079.0833
$ws.Range("N40").Value = -9634
$ws.Range("H46").Value = 4118.7334
$ws.Range("I46").Value = 933.5
$ws.Range("J46").Value = 6242.222
$ws.Range("K46").Value = 933.5
$ws.Range("L46").Value = 6242.222
$ws.Range("M46").Value = -745.5
$ws.Range("N46").Value = -6618.222
$ws.Range("H82").Value = 1901.8077
$ws.Range("I82").Value = 1830.8096
$ws.Range("J82").Value = 2200
$ws.Range("K82").Value = 1830.8096
$ws.Range("L82").Value = 2200
$ws.Range("M82").Value = -1469.8096
$ws.Range("N82").Value = -2922
$ws.Range("H85").Value = 1901.8077
$ws.Range("I85").Value = 1830.8096
$ws.Range("J85").Value = 2200
$ws.Range("K85").Value = 1830.8096
$ws.Range("L85").Value = 2200
$ws.Range("M85").Value = -582.8096
$ws.Range("N85").Value = -4696
$ws.Range("H126").Value = 2986.6667
$ws.Range("I126").Value = 2895
$ws.Range("J126").Value = 3170
$ws.Range("K126").Value = 8685
$ws.Range("L126").Value = 9510
$ws.Range("M126").Value = -6215
$ws.Range("N126").Value = -14450
$ws.Range("H136").Value = 1954.6428
$ws.Range("I136").Value = 1556
$ws.Range("K136").Value = 4668
$ws.Range("M136").Value = -2118

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 916.7143
$ws.Range("I81").Value = 902.8333
$ws.Range("K81").Value = 1805.6666
$ws.Range("M81").Value = -744.6666
$ws.Range("H84").Value = 916.7143
$ws.Range("I84").Value = 902.8333
$ws.Range("K84").Value = 9028.333000000001
$ws.Range("M84").Value = -3724.333000000001
$ws.Range("H126").Value = 50000790
$ws.Range("I126").Value = 58824284
$ws.Range("J126").Value = 1001.6667
$ws.Range("K126").Value = 176472852
$ws.Range("L126").Value = 3005.0001
$ws.Range("M126").Value = -176470382
$ws.Range("N126").Value = -7945.0001
$ws.Range("H136").Value = 1135.9584
$ws.Range("I136").Value = 793.15
$ws.Range("J136").Value = 2850
$ws.Range("K136").Value = 2379.45
$ws.Range("L136").Value = 8550
$ws.Range("M136").Value = 170.5500000000002
$ws.Range("N136").Value = -13650
